$wb = $excel.ActiveWorkbook

# ALC!row28
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 16730.77  # H28: was 12920.8125
$ws.Cells.Item(28, 9).Value = 2832  # I28: was 2361.6365
$ws.Cells.Item(28, 10).Value = 48003  # J28: was 36151
$ws.Cells.Item(28, 11).Value = 2832  # K28: was 2361.6365
$ws.Cells.Item(28, 12).Value = 48003  # L28: was 36151
$ws.Cells.Item(28, 13).Value = -2347  # M28: was -1876.6365
$ws.Cells.Item(28, 14).Value = -48973  # N28: was -37121

# ALC!row70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1000.6667  # H70: was 755.55554
$ws.Cells.Item(70, 9).Value = 656.44446  # I70: was 525
$ws.Cells.Item(70, 10).Value = 2033.3334  # J70: was 2600
$ws.Cells.Item(70, 11).Value = 1969.33338  # K70: was 1575
$ws.Cells.Item(70, 12).Value = 6100.0002  # L70: was 7800
$ws.Cells.Item(70, 13).Value = -1699.33338  # M70: was -1305
$ws.Cells.Item(70, 14).Value = -6640.0002  # N70: was -8340

# ALC!row73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 1000.6667  # H73: was 755.55554
$ws.Cells.Item(73, 9).Value = 656.44446  # I73: was 525
$ws.Cells.Item(73, 10).Value = 2033.3334  # J73: was 2600
$ws.Cells.Item(73, 11).Value = 1969.33338  # K73: was 1575
$ws.Cells.Item(73, 12).Value = 6100.0002  # L73: was 7800
$ws.Cells.Item(73, 13).Value = -1033.33338  # M73: was -639
$ws.Cells.Item(73, 14).Value = -7972.0002  # N73: was -9672

# ALC!row100
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 4498.75  # H100: was 4562.5
$ws.Cells.Item(100, 9).Value = 2798  # I100: was 2888.889
$ws.Cells.Item(100, 10).Value = 7333.3335  # J100: was 6714.2856
$ws.Cells.Item(100, 11).Value = 2798  # K100: was 2888.889
$ws.Cells.Item(100, 12).Value = 7333.3335  # L100: was 6714.2856
$ws.Cells.Item(100, 13).Value = -2257  # M100: was -2347.889
$ws.Cells.Item(100, 14).Value = -8415.333500000001  # N100: was -7796.2856

# ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 943.3409  # H129: was 994.5111000000001
$ws.Cells.Item(129, 9).Value = 423.64285  # I129: was 474.27274
$ws.Cells.Item(129, 10).Value = 1185.8667  # J129: was 1162.8235
$ws.Cells.Item(129, 11).Value = 1270.92855  # K129: was 1422.81822
$ws.Cells.Item(129, 12).Value = 3557.6001  # L129: was 3488.4705
$ws.Cells.Item(129, 13).Value = 3729.07145  # M129: was 3577.18178
$ws.Cells.Item(129, 14).Value = -13557.6001  # N129: was -13488.4705

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1226.2712  # H137: was 1344.6227
$ws.Cells.Item(137, 9).Value = 1018.75  # I137: was 1186.3334
$ws.Cells.Item(137, 11).Value = 3056.25  # K137: was 3559.0002
$ws.Cells.Item(137, 13).Value = -506.25  # M137: was -1009.0002

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5692005.5  # H32: was 2604.12
$ws.Cells.Item(32, 9).Value = 6587328  # I32: was 2604.12
$ws.Cells.Item(32, 10).Value = 21630.584  # J32: was 0
$ws.Cells.Item(32, 11).Value = 6587328  # K32: was 2604.12
$ws.Cells.Item(32, 12).Value = 21630.584  # L32: was 0
$ws.Cells.Item(32, 13).Value = -6587041  # M32: was -2317.12
$ws.Cells.Item(32, 14).Value = -22204.584  # N32: was None

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1820.6207  # H74: was 1907.7455
$ws.Cells.Item(74, 9).Value = 1290.1666  # I74: was 1336.9714
$ws.Cells.Item(74, 10).Value = 2688.6365  # J74: was 2906.6
$ws.Cells.Item(74, 11).Value = 1290.1666  # K74: was 1336.9714
$ws.Cells.Item(74, 12).Value = 2688.6365  # L74: was 2906.6
$ws.Cells.Item(74, 13).Value = -416.1666  # M74: was -462.9713999999999
$ws.Cells.Item(74, 14).Value = -4436.636500000001  # N74: was -4654.6

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1820.6207  # H77: was 1907.7455
$ws.Cells.Item(77, 9).Value = 1290.1666  # I77: was 1336.9714
$ws.Cells.Item(77, 10).Value = 2688.6365  # J77: was 2906.6
$ws.Cells.Item(77, 11).Value = 6450.833000000001  # K77: was 6684.857
$ws.Cells.Item(77, 12).Value = 13443.1825  # L77: was 14533
$ws.Cells.Item(77, 13).Value = -2082.833000000001  # M77: was -2316.857
$ws.Cells.Item(77, 14).Value = -22179.1825  # N77: was -23269

# ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 113220.75  # H110: was 82336.37
$ws.Cells.Item(110, 9).Value = 150644.33  # I110: was 82336.37
$ws.Cells.Item(110, 10).Value = 950  # J110: was 0
$ws.Cells.Item(110, 11).Value = 150644.33  # K110: was 82336.37
$ws.Cells.Item(110, 12).Value = 950  # L110: was 0
$ws.Cells.Item(110, 13).Value = -148599.33  # M110: was -80291.37
$ws.Cells.Item(110, 14).Value = -5040  # N110: was None

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 4482.4473  # H132: was 5683.893
$ws.Cells.Item(132, 9).Value = 4169.12  # I132: was 5374.3335
$ws.Cells.Item(132, 10).Value = 5085  # J132: was 6241.1
$ws.Cells.Item(132, 11).Value = 12507.36  # K132: was 16123.0005
$ws.Cells.Item(132, 12).Value = 15255  # L132: was 18723.3
$ws.Cells.Item(132, 13).Value = -9977.360000000001  # M132: was -13593.0005
$ws.Cells.Item(132, 14).Value = -20315  # N132: was -23783.3

# BSM!row64
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 619.36365  # H64: was 718.8333
$ws.Cells.Item(64, 9).Value = 588.25  # I64: was 676.5
$ws.Cells.Item(64, 10).Value = 702.3333  # J64: was 803.5
$ws.Cells.Item(64, 11).Value = 588.25  # K64: was 676.5
$ws.Cells.Item(64, 12).Value = 702.3333  # L64: was 803.5
$ws.Cells.Item(64, 13).Value = -363.25  # M64: was -451.5
$ws.Cells.Item(64, 14).Value = -1152.3333  # N64: was -1253.5

# BSM!row67
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(67, 8).Value = 619.36365  # H67: was 718.8333
$ws.Cells.Item(67, 9).Value = 588.25  # I67: was 676.5
$ws.Cells.Item(67, 10).Value = 702.3333  # J67: was 803.5
$ws.Cells.Item(67, 11).Value = 588.25  # K67: was 676.5
$ws.Cells.Item(67, 12).Value = 702.3333  # L67: was 803.5
$ws.Cells.Item(67, 13).Value = 191.75  # M67: was 103.5
$ws.Cells.Item(67, 14).Value = -2262.3333  # N67: was -2363.5

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 0  # H99: was 1312.3334
$ws.Cells.Item(99, 9).Value = 0  # I99: was 1171.4286
$ws.Cells.Item(99, 10).Value = 0  # J99: was 1805.5
$ws.Cells.Item(99, 11).Value = 0  # K99: was 1171.4286
$ws.Cells.Item(99, 12).Value = 0  # L99: was 1805.5
$ws.Cells.Item(99, 13).ClearContents()  # M99: was 326.5714
$ws.Cells.Item(99, 14).ClearContents()  # N99: was -4801.5

# CRP!row18
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(18, 8).Value = 39558  # H18: was 40000
$ws.Cells.Item(18, 10).Value = 39558  # J18: was 40000
$ws.Cells.Item(18, 12).Value = 39558  # L18: was 40000
$ws.Cells.Item(18, 14).Value = -40018  # N18: was -40460

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4116.314  # H31: was 4188.8115
$ws.Cells.Item(31, 9).Value = 1442.4722  # I31: was 1527.4242
$ws.Cells.Item(31, 10).Value = 6041.48  # J31: was 5877.769
$ws.Cells.Item(31, 11).Value = 1442.4722  # K31: was 1527.4242
$ws.Cells.Item(31, 12).Value = 6041.48  # L31: was 5877.769
$ws.Cells.Item(31, 13).Value = -1147.4722  # M31: was -1232.4242
$ws.Cells.Item(31, 14).Value = -6631.48  # N31: was -6467.769

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 4116.314  # H34: was 4188.8115
$ws.Cells.Item(34, 9).Value = 1442.4722  # I34: was 1527.4242
$ws.Cells.Item(34, 10).Value = 6041.48  # J34: was 5877.769
$ws.Cells.Item(34, 11).Value = 1442.4722  # K34: was 1527.4242
$ws.Cells.Item(34, 12).Value = 6041.48  # L34: was 5877.769
$ws.Cells.Item(34, 13).Value = -1240.4722  # M34: was -1325.4242
$ws.Cells.Item(34, 14).Value = -6445.48  # N34: was -6281.769

# CRP!row62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 4524.303  # H62: was 4331.467
$ws.Cells.Item(62, 9).Value = 4541.9033  # I62: was 4363.7393
$ws.Cells.Item(62, 10).Value = 4251.5  # J62: was 3960.3333
$ws.Cells.Item(62, 11).Value = 4541.9033  # K62: was 4363.7393
$ws.Cells.Item(62, 12).Value = 4251.5  # L62: was 3960.3333
$ws.Cells.Item(62, 13).Value = -3917.9033  # M62: was -3739.7393
$ws.Cells.Item(62, 14).Value = -5499.5  # N62: was -5208.3333

# CRP!row65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 4524.303  # H65: was 4331.467
$ws.Cells.Item(65, 9).Value = 4541.9033  # I65: was 4363.7393
$ws.Cells.Item(65, 10).Value = 4251.5  # J65: was 3960.3333
$ws.Cells.Item(65, 11).Value = 22709.5165  # K65: was 21818.6965
$ws.Cells.Item(65, 12).Value = 21257.5  # L65: was 19801.6665
$ws.Cells.Item(65, 13).Value = -19589.5165  # M65: was -18698.6965
$ws.Cells.Item(65, 14).Value = -27497.5  # N65: was -26041.6665

# CRP!row86
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 2392.04  # H86: was 2218
$ws.Cells.Item(86, 9).Value = 2468.5789  # I86: was 2210.1
$ws.Cells.Item(86, 10).Value = 2149.6667  # J86: was 2249.6
$ws.Cells.Item(86, 11).Value = 2468.5789  # K86: was 2210.1
$ws.Cells.Item(86, 12).Value = 2149.6667  # L86: was 2249.6
$ws.Cells.Item(86, 13).Value = -1345.5789  # M86: was -1087.1
$ws.Cells.Item(86, 14).Value = -4395.6667  # N86: was -4495.6

# CRP!row89
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 2392.04  # H89: was 2218
$ws.Cells.Item(89, 9).Value = 2468.5789  # I89: was 2210.1
$ws.Cells.Item(89, 10).Value = 2149.6667  # J89: was 2249.6
$ws.Cells.Item(89, 11).Value = 12342.8945  # K89: was 11050.5
$ws.Cells.Item(89, 12).Value = 10748.3335  # L89: was 11248
$ws.Cells.Item(89, 13).Value = -6726.8945  # M89: was -5434.5
$ws.Cells.Item(89, 14).Value = -21980.3335  # N89: was -22480

# CRP!row120
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(120, 8).Value = 30001  # H120: was 30000.5
$ws.Cells.Item(120, 10).Value = 0  # J120: was 30000
$ws.Cells.Item(120, 12).Value = 0  # L120: was 30000
$ws.Cells.Item(120, 14).ClearContents()  # N120: was -37258

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 3393.1064  # H131: was 3395.5625
$ws.Cells.Item(131, 9).Value = 513.0769  # I131: was 560.9091
$ws.Cells.Item(131, 10).Value = 4494.294  # J131: was 4238.2974
$ws.Cells.Item(131, 11).Value = 1539.2307  # K131: was 1682.7273
$ws.Cells.Item(131, 12).Value = 13482.882  # L131: was 12714.8922
$ws.Cells.Item(131, 13).Value = 3500.7693  # M131: was 3357.2727
$ws.Cells.Item(131, 14).Value = -23562.882  # N131: was -22794.8922

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1827  # H102: was 2022.55
$ws.Cells.Item(102, 9).Value = 1670.4736  # I102: was 1907.6666
$ws.Cells.Item(102, 10).Value = 2322.6667  # J102: was 2367.2
$ws.Cells.Item(102, 11).Value = 1670.4736  # K102: was 1907.6666
$ws.Cells.Item(102, 12).Value = 2322.6667  # L102: was 2367.2
$ws.Cells.Item(102, 13).Value = -48.47360000000003  # M102: was -285.6666
$ws.Cells.Item(102, 14).Value = -5566.6667  # N102: was -5611.2

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 6096.2915  # H122: was 7340.15
$ws.Cells.Item(122, 9).Value = 6300.5  # I122: was 8866.933999999999
$ws.Cells.Item(122, 10).Value = 3850  # J122: was 2759.8
$ws.Cells.Item(122, 11).Value = 18901.5  # K122: was 26600.802
$ws.Cells.Item(122, 12).Value = 11550  # L122: was 8279.400000000001
$ws.Cells.Item(122, 13).Value = -16451.5  # M122: was -24150.802
$ws.Cells.Item(122, 14).Value = -16450  # N122: was -13179.4

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2207.261  # H132: was 2075.4897
$ws.Cells.Item(132, 9).Value = 1724.2285  # I132: was 1635.3414
$ws.Cells.Item(132, 10).Value = 3744.182  # J132: was 4331.25
$ws.Cells.Item(132, 11).Value = 5172.6855  # K132: was 4906.0242
$ws.Cells.Item(132, 12).Value = 11232.546  # L132: was 12993.75
$ws.Cells.Item(132, 13).Value = -2642.6855  # M132: was -2376.0242
$ws.Cells.Item(132, 14).Value = -16292.546  # N132: was -18053.75

# LTW!row40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 1518.1428  # H40: was 1560.9166
$ws.Cells.Item(40, 9).Value = 1437.8334  # I40: was 1521
$ws.Cells.Item(40, 11).Value = 1437.8334  # K40: was 1521
$ws.Cells.Item(40, 13).Value = -1301.8334  # M40: was -1385

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2914.2856  # H46: was 2950
$ws.Cells.Item(46, 9).Value = 1000  # I46: was 675
$ws.Cells.Item(46, 10).Value = 4350  # J46: was 7500
$ws.Cells.Item(46, 11).Value = 1000  # K46: was 675
$ws.Cells.Item(46, 12).Value = 4350  # L46: was 7500
$ws.Cells.Item(46, 13).Value = -812  # M46: was -487
$ws.Cells.Item(46, 14).Value = -4726  # N46: was -7876

# LTW!row61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2525.5356  # H61: was 2831.7917
$ws.Cells.Item(61, 9).Value = 971.9231  # I61: was 1046.3
$ws.Cells.Item(61, 10).Value = 3872  # J61: was 4107.143
$ws.Cells.Item(61, 11).Value = 971.9231  # K61: was 1046.3
$ws.Cells.Item(61, 12).Value = 3872  # L61: was 4107.143
$ws.Cells.Item(61, 13).Value = -769.9231  # M61: was -844.3
$ws.Cells.Item(61, 14).Value = -4276  # N61: was -4511.143

# LTW!row100
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 47235.2  # H100: was 45014.953
$ws.Cells.Item(100, 9).Value = 51998.89  # I100: was 49294.21
$ws.Cells.Item(100, 11).Value = 51998.89  # K100: was 49294.21
$ws.Cells.Item(100, 13).Value = -51457.89  # M100: was -48753.21

# LTW!row113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 2525.5356  # H113: was 2831.7917
$ws.Cells.Item(113, 9).Value = 971.9231  # I113: was 1046.3
$ws.Cells.Item(113, 10).Value = 3872  # J113: was 4107.143
$ws.Cells.Item(113, 11).Value = 971.9231  # K113: was 1046.3
$ws.Cells.Item(113, 12).Value = 3872  # L113: was 4107.143
$ws.Cells.Item(113, 13).Value = 1198.0769  # M113: was 1123.7
$ws.Cells.Item(113, 14).Value = -8212  # N113: was -8447.143

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1494.4445  # H122: was 2089
$ws.Cells.Item(122, 9).Value = 1488.1428  # I122: was 2256.9285
$ws.Cells.Item(122, 10).Value = 1516.5  # J122: was 1618.8
$ws.Cells.Item(122, 11).Value = 4464.428400000001  # K122: was 6770.7855
$ws.Cells.Item(122, 12).Value = 4549.5  # L122: was 4856.4
$ws.Cells.Item(122, 13).Value = -2014.428400000001  # M122: was -4320.7855
$ws.Cells.Item(122, 14).Value = -9449.5  # N122: was -9756.4

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1544.909  # H126: was 1545.2727
$ws.Cells.Item(126, 9).Value = 1632.6666  # I126: was 1633.1111
$ws.Cells.Item(126, 11).Value = 4897.9998  # K126: was 4899.3333
$ws.Cells.Item(126, 13).Value = -2427.9998  # M126: was -2429.3333

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2734268.2  # H132: was 2734551
$ws.Cells.Item(132, 9).Value = 1537.4286  # I132: was 2014.3954
$ws.Cells.Item(132, 10).Value = 8775042  # J132: was 9262278
$ws.Cells.Item(132, 11).Value = 4612.2858  # K132: was 6043.1862
$ws.Cells.Item(132, 12).Value = 26325126  # L132: was 27786834
$ws.Cells.Item(132, 13).Value = -2082.2858  # M132: was -3513.1862
$ws.Cells.Item(132, 14).Value = -26330186  # N132: was -27791894
